# Task 2.2.8 - add Inflammatory Markers related concepts + refresh metadata
# (alcohol/caffeine/substance-use concepts plus the HRV-inflammation
# correlation concept, bumped Date/Count on the Metadata sheet).

$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item(1)       # "Metadata"
$concepts = $wb.Worksheets.Item(2)   # "Concepts"

# ---------------------------------------------------------------------
# 1. Metadata sheet: refresh the generation Date and the concept Count.
# ---------------------------------------------------------------------
$meta.Cells.Item(8, 2).Value = "2025-11-27T11:57:11+00:00"

# "Count" (B27) holds a numeric-looking string ("54" -> "62") that must stay
# a text cell (it always has been). Force text via a quote-prefixed value,
# then restore the row's normal (non quote-prefixed) look by re-applying the
# formatting already used throughout the sheet, so the cell style itself
# isn't left pointing at a brand-new "quote prefix" style.
$meta.Cells.Item(27, 2).Value = "'62"
$meta.Cells.Item(26, 2).Copy()
$meta.Cells.Item(27, 2).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Concepts sheet: append the 8 new rows (56-63).
# ---------------------------------------------------------------------
$newRows = @(
    @("alcohol-type", "Type of alcoholic beverage", "Type of alcoholic beverage consumed (beer, wine, spirits, etc.)"),
    @("caffeine-source", "Caffeine source", "Source of dietary caffeine intake (coffee, tea, energy drinks, etc.)"),
    @("last-caffeine-time", "Time of last caffeine intake", "Date and time of the most recent caffeine consumption"),
    @("substance-frequency", "Frequency of substance use", "How often a substance is used (daily, weekly, monthly, etc.)"),
    @("last-use-date", "Date of last substance use", "Date when a substance was last used"),
    @("substance-use-summary", "Substance use summary", "Comprehensive summary of all substance use patterns for lifestyle medicine assessment"),
    @("substance-risk-level", "Overall substance use risk level", "Aggregate risk assessment based on all substance use patterns"),
    @("hrv-inflammation-correlation", "HRV-Inflammation correlation assessment", "Assessment of the correlation between heart rate variability metrics (especially RMSSD) and inflammatory biomarkers (CRP, IL-6). Based on RS1 systematic review finding: inverse correlation between vagal tone and systemic inflammation via cholinergic anti-inflammatory pathway (Tracey 2002).")
)

$templateRow = 55
$destRow = $templateRow

foreach ($row in $newRows) {
    $destRow = $destRow + 1

    $srcRange = $concepts.Range("A$templateRow" + ":D$templateRow")
    $dstRange = $concepts.Range("A$destRow" + ":D$destRow")

    # Copy the "Level" value (always the shared string "1") together with
    # the row's style in two passes, exactly like the existing rows, so the
    # new row's cells keep t="s" (instead of turning into a Number) without
    # growing the style table.
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4163)  # xlPasteValues
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = $false

    $concepts.Cells.Item($destRow, 2).Value = $row[0]
    $concepts.Cells.Item($destRow, 3).Value = $row[1]
    $concepts.Cells.Item($destRow, 4).Value = $row[2]
}
